$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1953.1818
$ws.Range("I41").Value = 2378.1428
$ws.Range("J41").Value = 1209.5
$ws.Range("K41").Value = 2378.1428
$ws.Range("L41").Value = 1209.5
$ws.Range("M41").Value = -1938.1428
$ws.Range("N41").Value = -2089.5
$ws.Range("H62").Value = 4001
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 4001
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H69").Value = 12668.066
$ws.Range("J69").Value = 13286.286
$ws.Range("L69").Value = 39858.858
$ws.Range("N69").Value = -41606.858
$ws.Range("H72").Value = 12668.066
$ws.Range("J72").Value = 13286.286
$ws.Range("L72").Value = 119576.574
$ws.Range("N72").Value = -128312.574
$ws.Range("H74").Value = 18525594
$ws.Range("J74").Value = 7175
$ws.Range("L74").Value = 7175
$ws.Range("N74").Value = -9047
$ws.Range("H77").Value = 18525594
$ws.Range("J77").Value = 7175
$ws.Range("L77").Value = 35875
$ws.Range("N77").Value = -45235
$ws.Range("H97").Value = 10999
$ws.Range("J97").Value = 10999
$ws.Range("L97").Value = 32997
$ws.Range("N97").Value = -33989
$ws.Range("H132").Value = 2437.4285
$ws.Range("I132").Value = 2460.2927
$ws.Range("K132").Value = 7380.8781
$ws.Range("M132").Value = -4850.8781
$ws.Range("H137").Value = 5533.4585
$ws.Range("I137").Value = 5165.15
$ws.Range("K137").Value = 15495.45
$ws.Range("M137").Value = -12945.45
$ws.Range("H138").Value = 6744.3135
$ws.Range("I138").Value = 6386.3335
$ws.Range("J138").Value = 6761.0938
$ws.Range("K138").Value = 19159.0005
$ws.Range("L138").Value = 20283.2814
$ws.Range("M138").Value = -14019.0005
$ws.Range("N138").Value = -30563.2814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21848.783
$ws.Range("I32").Value = 8365.825999999999
$ws.Range("J32").Value = 43999.355
$ws.Range("K32").Value = 8365.825999999999
$ws.Range("L32").Value = 43999.355
$ws.Range("M32").Value = -8078.825999999999
$ws.Range("N32").Value = -44573.355
$ws.Range("H61").Value = 41671776
$ws.Range("I61").Value = 45459668
$ws.Range("K61").Value = 45459668
$ws.Range("M61").Value = -45459456
$ws.Range("H92").Value = 63333
$ws.Range("J92").Value = 63333
$ws.Range("L92").Value = 63333
$ws.Range("N92").Value = -68325
$ws.Range("H136").Value = 41671776
$ws.Range("I136").Value = 45459668
$ws.Range("K136").Value = 136379004
$ws.Range("M136").Value = -136376454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6732.72
$ws.Range("I58").Value = 7238.143
$ws.Range("K58").Value = 7238.143
$ws.Range("M58").Value = -7035.143
$ws.Range("H99").Value = 4789.231
$ws.Range("I99").Value = 4252.4443
$ws.Range("K99").Value = 4252.4443
$ws.Range("M99").Value = -2754.4443
$ws.Range("H107").Value = 612737.5600000001
$ws.Range("I107").Value = 1020883.5
$ws.Range("K107").Value = 1020883.5
$ws.Range("M107").Value = -1018963.5
$ws.Range("H122").Value = 4518.295
$ws.Range("I122").Value = 4533.1924
$ws.Range("J122").Value = 4432.222
$ws.Range("K122").Value = 13599.5772
$ws.Range("L122").Value = 13296.666
$ws.Range("M122").Value = -11149.5772
$ws.Range("N122").Value = -18196.666
$ws.Range("H126").Value = 4789.231
$ws.Range("I126").Value = 4252.4443
$ws.Range("K126").Value = 12757.3329
$ws.Range("M126").Value = -10287.3329
$ws.Range("H132").Value = 1921.2623
$ws.Range("I132").Value = 1790.9814
$ws.Range("J132").Value = 2926.2856
$ws.Range("K132").Value = 5372.9442
$ws.Range("L132").Value = 8778.856800000001
$ws.Range("M132").Value = -2842.9442
$ws.Range("N132").Value = -13838.8568
$ws.Range("H134").Value = 2372.8484
$ws.Range("I134").Value = 2452.724
$ws.Range("K134").Value = 7358.172
$ws.Range("M134").Value = -4823.172
$ws.Range("H136").Value = 6732.72
$ws.Range("I136").Value = 7238.143
$ws.Range("K136").Value = 21714.429
$ws.Range("M136").Value = -19164.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1341
$ws.Range("I118").Value = 2474.5
$ws.Range("J118").Value = 207.5
$ws.Range("K118").Value = 7423.5
$ws.Range("L118").Value = 622.5
$ws.Range("M118").Value = -6180.5
$ws.Range("N118").Value = -3108.5
$ws.Range("H129").Value = 4110.3335
$ws.Range("I129").Value = 2899
$ws.Range("J129").Value = 4261.75
$ws.Range("K129").Value = 8697
$ws.Range("L129").Value = 12785.25
$ws.Range("M129").Value = -3697
$ws.Range("N129").Value = -22785.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1036.7222
$ws.Range("I102").Value = 1033.75
$ws.Range("K102").Value = 1033.75
$ws.Range("M102").Value = 588.25
$ws.Range("H113").Value = 5474.067
$ws.Range("I113").Value = 4460.2856
$ws.Range("J113").Value = 6361.125
$ws.Range("K113").Value = 4460.2856
$ws.Range("L113").Value = 6361.125
$ws.Range("M113").Value = -2290.2856
$ws.Range("N113").Value = -10701.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7604.4
$ws.Range("J7").Value = 9591.571
$ws.Range("L7").Value = 9591.571
$ws.Range("N7").Value = -9815.571
$ws.Range("H22").Value = 1898.2
$ws.Range("I22").Value = 1209.3529
$ws.Range("K22").Value = 1209.3529
$ws.Range("M22").Value = -914.3529000000001
$ws.Range("H27").Value = 1898.2
$ws.Range("I27").Value = 1209.3529
$ws.Range("K27").Value = 1209.3529
$ws.Range("M27").Value = -1102.3529
$ws.Range("H46").Value = 1705.9231
$ws.Range("I46").Value = 1562.6364
$ws.Range("K46").Value = 1562.6364
$ws.Range("M46").Value = -1374.6364
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H55").Value = 621.10345
$ws.Range("J55").Value = 853.2143
$ws.Range("L55").Value = 853.2143
$ws.Range("N55").Value = -1199.2143
$ws.Range("H61").Value = 4930.5938
$ws.Range("I61").Value = 4814.0356
$ws.Range("J61").Value = 5746.5
$ws.Range("K61").Value = 4814.0356
$ws.Range("L61").Value = 5746.5
$ws.Range("M61").Value = -4612.0356
$ws.Range("N61").Value = -6150.5
$ws.Range("H98").Value = 28855
$ws.Range("J98").Value = 28855
$ws.Range("L98").Value = 28855
$ws.Range("N98").Value = -34845
$ws.Range("H100").Value = 6398811
$ws.Range("I100").Value = 7346425.5
$ws.Range("J100").Value = 2410.75
$ws.Range("K100").Value = 7346425.5
$ws.Range("L100").Value = 2410.75
$ws.Range("M100").Value = -7345884.5
$ws.Range("N100").Value = -3492.75
$ws.Range("H113").Value = 4930.5938
$ws.Range("I113").Value = 4814.0356
$ws.Range("J113").Value = 5746.5
$ws.Range("K113").Value = 4814.0356
$ws.Range("L113").Value = 5746.5
$ws.Range("M113").Value = -2644.0356
$ws.Range("N113").Value = -10086.5
$ws.Range("H126").Value = 7604.4
$ws.Range("J126").Value = 9591.571
$ws.Range("L126").Value = 28774.713
$ws.Range("N126").Value = -33714.713
$ws.Range("H132").Value = 20008168
$ws.Range("I132").Value = 21913184
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 65739552
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -65737022
$ws.Range("N132").Value = -21560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5502.3096
$ws.Range("I132").Value = 3225.743
$ws.Range("J132").Value = 16885.143
$ws.Range("K132").Value = 9677.228999999999
$ws.Range("L132").Value = 50655.429
$ws.Range("M132").Value = -7147.228999999999
$ws.Range("N132").Value = -55715.429
$ws.Range("H136").Value = 2256.3809
$ws.Range("J136").Value = 2550
$ws.Range("L136").Value = 7650
$ws.Range("N136").Value = -12750
